$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 8 gets promoted to the "bordered" style (same look as rows 4-6):
#    copy the formatting (borders/font) from row 4 onto row 8, keeping
#    row 8's own values/text in place.
# ---------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) New row 9 (plain / unbordered style, like rows 3 & 7).
#    NOTE: the literal two characters backslash+n (NOT a real line break)
#    are part of the source text, matching the game-script convention
#    used throughout this sheet.
# ---------------------------------------------------------------------
$ws.Range("B7:E7").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)

$ws.Range("B9").Value = 137

# ---------------------------------------------------------------------
# 3) New row 10 (bordered style, like row 4).
# ---------------------------------------------------------------------
$ws.Range("A4:E4").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

$ws.Range("B10").Value = 140

# ---------------------------------------------------------------------
# Shared-string insertion order matters (new strings are appended to
# sharedStrings.xml in the order they are first written), so write the
# English / Russian / "converted" columns column-by-column across both
# rows: C9, C10, D9, D10, E9, E10.
# ---------------------------------------------------------------------
$ws.Range("C9").Value = " I congratulate you on your guild\ngraduation."
$ws.Range("C10").Value = " I do hope you continue to rake in\nmoney...[K]and allow me to hoard it!"
$ws.Range("D9").Value = " Поздравляю вас с выпуском из\nгильдии."
$ws.Range("D10").Value = " Надеюсь, вы будете и дальше\nнести сюда деньги...[K] И позволите мне их\nхранить!"
$ws.Range("E9").Value = " Ðïèäñàâìÿý âàò ò âúðôòëïí éè\nãéìûäéé."
$ws.Range("E10").Value = " Îàäåýòû, âú áôäåóå é äàìûšå\nîåòóé òýäà äåîûãé...[K] É ðïèâïìéóå íîå éö\nöñàîéóû!"

$ws.Rows.Item(9).RowHeight = 21.6
$ws.Rows.Item(10).RowHeight = 31.8

# ---------------------------------------------------------------------
# 4) Selection moves to C7.
# ---------------------------------------------------------------------
$ws.Range("C7").Select()
